# Update C24 to append the Weekly Health Report reference, and move the
# active selection from C16 to C30, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C24").Value = "EEJ1_nagios_installation_and_montoring, EEJ1_Weekly_Health_Report"

$ws.Range("C30").Select()
